$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.497.70"
$ws.Range("E2").Value = "  -3.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.468.79"
$ws.Range("E3").Value = "  -2.32%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.07"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.42"
$ws.Range("E6").Value = "  -8.73%  "

$ws.Range("E7").Value = "  +2.90%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -3.21%  "

$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.05"
$ws.Range("E11").Value = "  -6.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000268"
$ws.Range("E12").Value = "  -2.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.04"
$ws.Range("E13").Value = "  -5.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.020.26"
$ws.Range("E14").Value = "  -2.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.470.95"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.06"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.451.13"
$ws.Range("E18").Value = "  -3.45%  "

$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.983"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "409.57"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("E22").Value = "  +2.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.54"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.15"
$ws.Range("E24").Value = "  -2.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.79"
$ws.Range("E26").Value = "  -5.90%  "

$ws.Range("E27").Value = "  -6.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.83"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "29.88"
$ws.Range("E29").Value = "  -2.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "613.76"
$ws.Range("E30").Value = "  -11.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.24"
$ws.Range("E31").Value = "  -8.99%  "

$ws.Range("E32").Value = "  -2.96%  "

$ws.Range("E33").Value = "  -4.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.52"
$ws.Range("E34").Value = "  -3.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.150"
$ws.Range("E35").Value = "  +6.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0779"
$ws.Range("E37").Value = "  -6.61%  "

# --- Rows 38/39 swap: Maker <-> InjectiveProtocol, with updated values ---
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.41"
$ws.Range("E38").Value = "  -7.10%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.315.90"
$ws.Range("E39").Value = "  +8.09%  "

$ws.Range("E40").Value = "  -6.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  -5.03%  "

$ws.Range("E44").Value = "  -2.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0410"
$ws.Range("E45").Value = "  -3.18%  "

$ws.Range("E46").Value = "  -8.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -1.78%  "

$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.24"
$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.30"
$ws.Range("E50").Value = "  -8.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.78"
$ws.Range("E51").Value = "  +5.73%  "
